$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.291.37"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "2.994.58"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.31"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.91"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.430"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.10"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.366"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "3.506.75"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.41"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "56.224.91"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "2.993.47"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.96"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.88"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.00"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.08"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.496"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").Value = "3.120.53"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "0.0₃0935"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.32"
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.27"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.86"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("E35").Value = "  -2.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.26"
$ws.Range("E37").Value = "  +7.54%  "
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "3.030.71"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.68"
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "2.188.48"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.42"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0851"
$ws.Range("E51").Value = "  -1.95%  "
